$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '328.69'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '1.40%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '41.22'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '4.98%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.627'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-1.48%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08200'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '2.35%'
$ws.Range('B6').Value = 'KuCoinToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '8.744'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.46%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '2.005'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.71%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '4.486'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-0.30%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.987'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '1.19%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9207'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-0.80%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1282'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '3.44%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.1954'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-1.01%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09229'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.78%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03889'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '8.15%'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.84%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001306'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.73%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.006320'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.56%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.448'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3485'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '0.42%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.242'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-5.34%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.1372'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '0.09%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.2412'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.02%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.04387'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.55%'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-0.45%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004314'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-6.74%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001201'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '4.34%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02780'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '11.00%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05405'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.89%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007697'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '0.85%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.008960'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-6.74%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002171'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '2.55%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01148'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '11.36%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006578'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-2.36%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.06%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.003211'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '8.15%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002282'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.40%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.06%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.06%'

Write-Output "Applied 74 cell updates"
